$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 44 (shifts the existing rows 44-72 down to 45-73)
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row with its values (same record as the old
# row 44, but with an updated Fecha/Volumen/Precios/Precio-$-Kg reading)
$ws.Range("A44").Value = 10
$ws.Range("B44").Value = "Vega Modelo de Temuco"
$ws.Range("C44").Value = "La Araucanía"
$ws.Range("D44").Value = 44606
$ws.Range("E44").Value = 9
$ws.Range("F44").Value = 100112030
$ws.Range("G44").Value = "Poroto granado"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 50
$ws.Range("K44").Value = 25000
$ws.Range("L44").Value = 25000
$ws.Range("M44").Value = 25000
$ws.Range("N44").Value = "$/saco 25 kilos"
$ws.Range("O44").Value = "Región de La Araucanía"
$ws.Range("P44").Value = 1000
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
